$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Columns B and F hold plain text that Excel's input parser will not
# mistake for a date/number, so they can be assigned directly.
$ws.Cells.Item($row, 2).Value = "10 Jan -- 16 Jan 2021"
$ws.Cells.Item($row, 6).Value = "KNN"

# Column A holds a date-formatted string ("2021-01-09"). Assigning it
# directly via .Value would make Excel auto-convert it into a date
# serial number (and would also stamp a new number-format style onto
# the cell). To keep it as plain text - exactly like the existing
# cells in this column - build it with a formula in a scratch cell,
# then copy/paste the computed value (not the formula, not formats)
# into the target cell.
$helper = $ws.Range("Z1")
$helper.Formula = '="2021-01-09"'
$helper.Copy()
$ws.Range("A49").PasteSpecial(-4163, -4142, $false, $false)
$helper.ClearContents()
$excel.CutCopyMode = $false

# Numeric columns.
$ws.Cells.Item($row, 3).Value = 3333.57
$ws.Cells.Item($row, 4).Value = 2251.19
$ws.Cells.Item($row, 5).Value = 1082.38
$ws.Cells.Item($row, 10).Value = 823.1
$ws.Cells.Item($row, 11).Value = 26.17
